$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 700
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("H6").Value = 1198.8
$ws.Range("I6").Value = 1198.8
$ws.Range("K6").Value = 3596.4
$ws.Range("M6").Value = -3484.4
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -1138
$ws.Range("H51").Value = 1444.4445
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -516
$ws.Range("H62").Value = 4024.25
$ws.Range("I62").Value = 2699
$ws.Range("K62").Value = 2699
$ws.Range("M62").Value = -2075
$ws.Range("H64").Value = 5999.8
$ws.Range("I64").Value = 5999.8
$ws.Range("K64").Value = 5999.8
$ws.Range("M64").Value = -5751.8
$ws.Range("H65").Value = 4024.25
$ws.Range("I65").Value = 2699
$ws.Range("K65").Value = 13495
$ws.Range("M65").Value = -10375
$ws.Range("H67").Value = 5999.8
$ws.Range("I67").Value = 5999.8
$ws.Range("K67").Value = 5999.8
$ws.Range("M67").Value = -5141.8
$ws.Range("H100").Value = 2525.8
$ws.Range("J100").Value = 1000
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082
$ws.Range("H129").Value = 1083
$ws.Range("I129").Value = 1083
$ws.Range("K129").Value = 3249
$ws.Range("M129").Value = 1751
$ws.Range("H132").Value = 1289
$ws.Range("I132").Value = 1289
$ws.Range("K132").Value = 3867
$ws.Range("M132").Value = -1337
$ws.Range("H135").Value = 1282.3636
$ws.Range("I135").Value = 456.33334
$ws.Range("K135").Value = 4107.00006
$ws.Range("M135").Value = -1572.00006
$ws.Range("H137").Value = 3513.6
$ws.Range("I137").Value = 3344.077
$ws.Range("J137").Value = 3828.4285
$ws.Range("K137").Value = 10032.231
$ws.Range("L137").Value = 11485.2855
$ws.Range("M137").Value = -7482.231
$ws.Range("N137").Value = -16585.2855
$ws.Range("H138").Value = 3942.1052
$ws.Range("I138").Value = 1081.5
$ws.Range("J138").Value = 5262.385
$ws.Range("K138").Value = 3244.5
$ws.Range("L138").Value = 15787.155
$ws.Range("M138").Value = 1895.5
$ws.Range("N138").Value = -26067.155
$ws.Range("H141").Value = 800
$ws.Range("I141").Value = 800
$ws.Range("K141").Value = 2400
$ws.Range("M141").Value = 2780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13110.792
$ws.Range("I32").Value = 13110.792
$ws.Range("K32").Value = 13110.792
$ws.Range("M32").Value = -12823.792
$ws.Range("H45").Value = 3022.5
$ws.Range("I45").Value = 3340.625
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 3340.625
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -2963.625
$ws.Range("N45").Value = -2504
$ws.Range("H61").Value = 1259.5
$ws.Range("I61").Value = 1349.375
$ws.Range("K61").Value = 1349.375
$ws.Range("M61").Value = -1137.375
$ws.Range("H74").Value = 18430.367
$ws.Range("I74").Value = 17885.666
$ws.Range("K74").Value = 17885.666
$ws.Range("M74").Value = -17011.666
$ws.Range("H77").Value = 18430.367
$ws.Range("I77").Value = 17885.666
$ws.Range("K77").Value = 89428.33
$ws.Range("M77").Value = -85060.33
$ws.Range("H97").Value = 125.42857
$ws.Range("I97").Value = 128.6
$ws.Range("J97").Value = 117.5
$ws.Range("K97").Value = 128.6
$ws.Range("L97").Value = 117.5
$ws.Range("M97").Value = 367.4
$ws.Range("N97").Value = -1109.5
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102774
$ws.Range("H136").Value = 1259.5
$ws.Range("I136").Value = 1349.375
$ws.Range("K136").Value = 4048.125
$ws.Range("M136").Value = -1498.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2221.889
$ws.Range("I31").Value = 2059.6
$ws.Range("K31").Value = 2059.6
$ws.Range("M31").Value = -1764.6
$ws.Range("H34").Value = 2221.889
$ws.Range("I34").Value = 2059.6
$ws.Range("K34").Value = 2059.6
$ws.Range("M34").Value = -1857.6
$ws.Range("H58").Value = 5042.727
$ws.Range("I58").Value = 4847
$ws.Range("K58").Value = 4847
$ws.Range("M58").Value = -4644
$ws.Range("H86").Value = 7349.25
$ws.Range("I86").Value = 6499.6665
$ws.Range("J86").Value = 9898
$ws.Range("K86").Value = 6499.6665
$ws.Range("L86").Value = 9898
$ws.Range("M86").Value = -5376.6665
$ws.Range("N86").Value = -12144
$ws.Range("H89").Value = 7349.25
$ws.Range("I89").Value = 6499.6665
$ws.Range("J89").Value = 9898
$ws.Range("K89").Value = 32498.3325
$ws.Range("L89").Value = 49490
$ws.Range("M89").Value = -26882.3325
$ws.Range("N89").Value = -60722
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = ""
$ws.Range("H106").Value = 90000
$ws.Range("J106").Value = 90000
$ws.Range("L106").Value = 90000
$ws.Range("N106").Value = -92524
$ws.Range("H132").Value = 1642.2858
$ws.Range("I132").Value = 1249.5
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 3748.5
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -1218.5
$ws.Range("N132").Value = -17057
$ws.Range("H134").Value = 6992
$ws.Range("I134").Value = 7990
$ws.Range("K134").Value = 23970
$ws.Range("M134").Value = -21435
$ws.Range("H136").Value = 5042.727
$ws.Range("I136").Value = 4847
$ws.Range("K136").Value = 14541
$ws.Range("M136").Value = -11991

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 587
$ws.Range("I14").Value = 587
$ws.Range("K14").Value = 1761
$ws.Range("M14").Value = -1588
$ws.Range("H15").Value = 299
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = 498
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 1494
$ws.Range("M15").Value = -160
$ws.Range("N15").Value = -1774
$ws.Range("H17").Value = 1680.6923
$ws.Range("I17").Value = 145
$ws.Range("J17").Value = 2363.2222
$ws.Range("K17").Value = 435
$ws.Range("L17").Value = 7089.6666
$ws.Range("M17").Value = -266
$ws.Range("N17").Value = -7427.6666
$ws.Range("H98").Value = 4123
$ws.Range("I98").Value = 4199.3335
$ws.Range("J98").Value = 3894
$ws.Range("K98").Value = 12598.0005
$ws.Range("L98").Value = 11682
$ws.Range("M98").Value = -11100.0005
$ws.Range("N98").Value = -14678
$ws.Range("H118").Value = 4593.4546
$ws.Range("I118").Value = 528.5
$ws.Range("J118").Value = 4999.95
$ws.Range("K118").Value = 1585.5
$ws.Range("L118").Value = 14999.85
$ws.Range("M118").Value = -342.5
$ws.Range("N118").Value = -17485.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4857.8335
$ws.Range("I113").Value = 5079.6
$ws.Range("K113").Value = 5079.6
$ws.Range("M113").Value = -2909.6
$ws.Range("H122").Value = 5663
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5521.8945
$ws.Range("I22").Value = 3942.3635
$ws.Range("J22").Value = 7693.75
$ws.Range("K22").Value = 3942.3635
$ws.Range("L22").Value = 7693.75
$ws.Range("M22").Value = -3647.3635
$ws.Range("N22").Value = -8283.75
$ws.Range("H27").Value = 5521.8945
$ws.Range("I27").Value = 3942.3635
$ws.Range("J27").Value = 7693.75
$ws.Range("K27").Value = 3942.3635
$ws.Range("L27").Value = 7693.75
$ws.Range("M27").Value = -3835.3635
$ws.Range("N27").Value = -7907.75
$ws.Range("H40").Value = 5500
$ws.Range("I40").Value = 5500
$ws.Range("K40").Value = 5500
$ws.Range("M40").Value = -5364
$ws.Range("H46").Value = 2697.182
$ws.Range("I46").Value = 1630
$ws.Range("K46").Value = 1630
$ws.Range("M46").Value = -1442

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("H122").Value = 4793.2856
$ws.Range("I122").Value = 4758.8335
$ws.Range("K122").Value = 14276.5005
$ws.Range("M122").Value = -11826.5005
$ws.Range("H132").Value = 1046.3334
$ws.Range("I132").Value = 1046.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3139.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -609.0001999999999
$ws.Range("N132").Value = ""
